# "NO and NaN preprocess added" - adds a second results sheet ("CompleteSVM")
# next to the existing "CompleteRFC2" sheet, holding the f1_score/accuracy
# results for the SVM run (mirrors the layout already used on CompleteRFC2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New worksheet, inserted right after the existing "CompleteRFC2" tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "CompleteSVM"

# Header row + result row (same A1:C2 layout as CompleteRFC2).
$ws2.Range("B1").Value = "f1_score"
$ws2.Range("C1").Value = "accuracy"
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = 0.4816025321822184
$ws2.Range("C2").Value = 0.7308

# Match the bold / centered / thin-bordered header style used on CompleteRFC2
# for the header cells (B1:C1) and the run-index cell (A2).
$headerRange = $ws2.Range("B1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$a2 = $ws2.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

# Keep the original sheet active/selected, as it was before the edit.
$ws1.Activate()
